# Updated cryptos list (Price / Volume(1h) columns) per the Sep 30 2024
# GitHub Actions data refresh. Cells D2:E51 hold plain display text (not
# numeric types) in the source workbook, so for any new value that Excel
# would otherwise auto-parse as a number we force the cell to Text via
# NumberFormat "@" before assigning it, preserving exact formatting
# (trailing zeros, thousands-dot groups, padded "%" strings, etc).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.842.07"
$ws.Range("E2").Value = "  -2.92%  "
$ws.Range("D3").Value = "2.621.77"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.14"
$ws.Range("E5").Value = "  -3.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.77"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.625"
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("E9").Value = "  -5.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.80"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.383"
$ws.Range("E11").Value = "  -3.75%  "
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.23"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").Value = "3.097.29"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000183"
$ws.Range("E15").Value = "  -6.90%  "
$ws.Range("D16").Value = "63.696.22"
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("D17").Value = "2.655.40"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.06"
$ws.Range("E18").Value = "  -3.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.66"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.56"
$ws.Range("E20").Value = "  -3.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "342.73"
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.37"
$ws.Range("E23").Value = "  -2.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.75"
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000109"
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "595.36"
$ws.Range("E26").Value = "  +5.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.19"
$ws.Range("E27").Value = "  -4.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.58"
$ws.Range("E28").Value = "  -1.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.161"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.87"
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.08"
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.73"
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.57"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.31"
$ws.Range("E35").Value = "  -2.40%  "
$ws.Range("E36").Value = "  -2.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.79"
$ws.Range("E37").Value = "  -3.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "153.98"
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.87"
$ws.Range("E40").Value = "  -3.52%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "41.46"
$ws.Range("E42").Value = "  -2.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.42"
$ws.Range("E43").Value = "  +6.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "156.33"
$ws.Range("E44").Value = "  -3.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.93"
$ws.Range("E45").Value = "  -3.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.89"
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0594"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.630"
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.101"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0249"
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.01"
$ws.Range("E51").Value = "  -3.81%  "
